# Fruta / hortaliza, semanal
# Insert a new data row at row 696 (pushing existing rows 696:758 down to 697:759)
# and populate it with the new week's price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(696).Insert()

$ws.Range("A696").Value = 6
$ws.Range("B696").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C696").Value = "Metropolitana"
$ws.Range("D696").Value = 45106
$ws.Range("E696").Value = 13
$ws.Range("F696").Value = 100112052
$ws.Range("G696").Value = "Albahaca"
$ws.Range("H696").Value = "Sin especificar"
$ws.Range("I696").Value = "Primera"
$ws.Range("J696").Value = 300
$ws.Range("K696").Value = 4500
$ws.Range("L696").Value = 5000
$ws.Range("M696").Value = 4733
$ws.Range("N696").Value = '$/paquete'
$ws.Range("O696").Value = "Región de Arica y Parinacota"
$ws.Range("P696").Value = 4733
$ws.Range("Q696").Value = 1
$ws.Range("R696").Value = "Hortaliza"
